# bug fix in Eduati data files
# SW1463_noCTRL_meas.xlsx: Sheet1 had stray leftover rows (45:87) containing
# only an incrementing index in column A beyond the real 44-row data table;
# remove them. Also the previously-active tab (Sheet3) becomes Sheet1, with
# the selection left where the user was working (near the bottom of the
# now-removed rows) before trimming.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Drop the stray rows 45:87 (only col A had data: a simple running index),
# shrinking the sheet's real extent back down to the 44-row data table.
$ws1.Range("A45:N87").EntireRow.Delete() | Out-Null

# Make Sheet1 the active/selected tab (previously Sheet3 was active), and
# leave the selection near where the trimmed rows used to be.
$ws1.Activate() | Out-Null
$ws1.Range("F65").Select() | Out-Null
